$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 30
$ws1.Range("F3").Value = 769
$ws1.Range("F4").Value = 31
$ws1.Range("F5").Value = 55
$ws1.Range("F6").Value = 49
$ws1.Range("F7").Value = 265
$ws1.Range("F8").Value = 3734
$ws1.Range("F10").Value = 4392
$ws1.Range("F11").Value = 488
$ws1.Range("F12").Value = 1085
$ws1.Range("F13").Value = 57

# Sheet "全部类型" (all types) - update "想去人数" (column F) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 30
$ws4.Range("F3").Value = 769
$ws4.Range("F4").Value = 31
$ws4.Range("F5").Value = 55
$ws4.Range("F6").Value = 49
$ws4.Range("F8").Value = 265
$ws4.Range("F9").Value = 3734
$ws4.Range("F11").Value = 4392
$ws4.Range("F12").Value = 488
$ws4.Range("F13").Value = 1085
$ws4.Range("F14").Value = 57
